$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text, preserving exact formatting (no
# locale-driven numeric/date auto-conversion, no leftover cell style diff).
function Set-TextValue {
    param($range, $text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.NumberFormat = "General"
    $range.Style = "Normal"
}

# --- Price / Volume(1h) updates for the refreshed crypto snapshot ---
Set-TextValue $ws.Range("D2") '61.931.25'
Set-TextValue $ws.Range("E2") '  -0.06%  '
Set-TextValue $ws.Range("D3") '2.418.32'
Set-TextValue $ws.Range("E3") '  +0.34%  '
Set-TextValue $ws.Range("E4") '  +0.08%  '
Set-TextValue $ws.Range("D5") '563.32'
Set-TextValue $ws.Range("E5") '  +1.30%  '
Set-TextValue $ws.Range("D6") '142.86'
Set-TextValue $ws.Range("E6") '  +0.43%  '
Set-TextValue $ws.Range("E7") '  -0.03%  '
Set-TextValue $ws.Range("D8") '0.531'
Set-TextValue $ws.Range("E8") '  +0.19%  '
Set-TextValue $ws.Range("E9") '  +0.73%  '
Set-TextValue $ws.Range("E10") '  -1.14%  '
Set-TextValue $ws.Range("D11") '5.21'
Set-TextValue $ws.Range("E11") '  -3.28%  '
Set-TextValue $ws.Range("D12") '0.350'
Set-TextValue $ws.Range("E12") '  -0.51%  '
Set-TextValue $ws.Range("D13") '25.81'
Set-TextValue $ws.Range("E13") '  -0.61%  '
Set-TextValue $ws.Range("D14") '0.0000172'
Set-TextValue $ws.Range("E14") '  -0.40%  '
Set-TextValue $ws.Range("D16") '61.860.42'
Set-TextValue $ws.Range("E16") '  +0.16%  '
Set-TextValue $ws.Range("D17") '2.411.51'
Set-TextValue $ws.Range("E17") '  +0.01%  '
Set-TextValue $ws.Range("D18") '11.30'
Set-TextValue $ws.Range("E18") '  +1.71%  '
Set-TextValue $ws.Range("D19") '322.99'
Set-TextValue $ws.Range("E19") '  +0.07%  '
Set-TextValue $ws.Range("E20") '  +1.59%  '
Set-TextValue $ws.Range("D21") '4.12'
Set-TextValue $ws.Range("E21") '  -1.20%  '
Set-TextValue $ws.Range("E22") '  +0.06%  '
Set-TextValue $ws.Range("D23") '66.52'
Set-TextValue $ws.Range("E23") '  +2.24%  '
Set-TextValue $ws.Range("D24") '1.73'
Set-TextValue $ws.Range("E24") '  +1.31%  '
Set-TextValue $ws.Range("D25") '8.79'
Set-TextValue $ws.Range("E25") '  -3.86%  '
Set-TextValue $ws.Range("D26") '556.80'
Set-TextValue $ws.Range("E26") '  -3.79%  '
Set-TextValue $ws.Range("E29") '  +1.10%  '
Set-TextValue $ws.Range("D30") '8.18'
Set-TextValue $ws.Range("E30") '  -0.94%  '
Set-TextValue $ws.Range("D31") '1.38'
Set-TextValue $ws.Range("E31") '  -4.31%  '
Set-TextValue $ws.Range("E32") '  -0.55%  '
Set-TextValue $ws.Range("E33") '  -0.20%  '
Set-TextValue $ws.Range("E34") '  -3.62%  '
Set-TextValue $ws.Range("D35") '1.00'
Set-TextValue $ws.Range("E35") '  -0.02%  '
Set-TextValue $ws.Range("D36") '4.74'
Set-TextValue $ws.Range("E36") '  +0.03%  '
Set-TextValue $ws.Range("E37") '  -1.02%  '
Set-TextValue $ws.Range("D38") '153.82'
Set-TextValue $ws.Range("E38") '  +1.80%  '
Set-TextValue $ws.Range("E39") '  -3.90%  '
Set-TextValue $ws.Range("D40") '18.51'
Set-TextValue $ws.Range("E40") '  -0.69%  '
Set-TextValue $ws.Range("D41") '1.81'
Set-TextValue $ws.Range("E41") '  -1.12%  '
Set-TextValue $ws.Range("D42") '0.993'
Set-TextValue $ws.Range("E42") '  -0.67%  '
Set-TextValue $ws.Range("E43") '  -2.12%  '
Set-TextValue $ws.Range("D44") '147.07'
Set-TextValue $ws.Range("E44") '  -2.11%  '
Set-TextValue $ws.Range("D45") '3.63'
Set-TextValue $ws.Range("E45") '  +0.01%  '
Set-TextValue $ws.Range("E46") '  -2.31%  '
Set-TextValue $ws.Range("E49") '  -0.16%  '
Set-TextValue $ws.Range("D50") '0.0227'
Set-TextValue $ws.Range("E50") '  -0.40%  '
Set-TextValue $ws.Range("D51") '11.58'
Set-TextValue $ws.Range("E51") '  +0.91%  '

# Row 27/28 swap (Binance-PegBSC-USD <-> WrappedeETH), with new D/E values
Set-TextValue $ws.Range("B27") 'WrappedeETH'
Set-TextValue $ws.Range("C27") 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
Set-TextValue $ws.Range("D27") '2.536.92'
Set-TextValue $ws.Range("E27") '  +0.35%  '

Set-TextValue $ws.Range("B28") 'Binance-PegBSC-USD'
Set-TextValue $ws.Range("C28") 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue $ws.Range("D28") '1.00'
Set-TextValue $ws.Range("E28") '  -0.04%  '

# Row 47/48 swap (Mantle <-> InjectiveProtocol), with new D/E values
Set-TextValue $ws.Range("B47") 'InjectiveProtocol'
Set-TextValue $ws.Range("C47") 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range("D47") '19.79'
Set-TextValue $ws.Range("E47") '  -1.68%  '

Set-TextValue $ws.Range("B48") 'Mantle'
Set-TextValue $ws.Range("C48") 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue $ws.Range("D48") '0.592'
Set-TextValue $ws.Range("E48") '  +0.73%  '

Write-Output "Done applying crypto list update."
